$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column H ("Productivity [MMGGE/yr]"), shifting I/J/K left to H/I/J
$ws.Columns.Item(8).Delete()

# Update the remaining numeric data cells with their new values
$ws.Range("C4").Value2 = -0.05487660350641402
$ws.Range("E4").Value2 = -0.02220258481033924
$ws.Range("F4").Value2 = -0.02655476221904888
$ws.Range("H4").Value2 = -0.01617875271501086
$ws.Range("I4").Value2 = -0.02424720098880395
$ws.Range("J4").Value2 = -0.05186680484024906

$ws.Range("C5").Value2 = -0.0378813675254701
$ws.Range("E5").Value2 = 0.00240375361501446
$ws.Range("F5").Value2 = -0.00284813939255757
$ws.Range("H5").Value2 = -0.03184649538598154
$ws.Range("I5").Value2 = 0.0001651206604826419
$ws.Range("J5").Value2 = 0.04313613676180762

$ws.Range("C6").Value2 = 0.03099045996183985
$ws.Range("E6").Value2 = -0.02410579242316969
$ws.Range("F6").Value2 = -0.02659844239376957
$ws.Range("H6").Value2 = -0.05781642326569306
$ws.Range("I6").Value2 = -0.02494521978087913
$ws.Range("J6").Value2 = 0.01808148117083027

$ws.Range("C7").Value2 = 0.04632191328765314
$ws.Range("E7").Value2 = 0.8729432997731991
$ws.Range("F7").Value2 = 0.8817736070944283
$ws.Range("H7").Value2 = 0.9999064956259824
$ws.Range("I7").Value2 = 0.8722964491857966
$ws.Range("J7").Value2 = 0.04654640618776307

$ws.Range("C8").Value2 = 0.97725856503426
$ws.Range("E8").Value2 = 0.01101998007992032
$ws.Range("F8").Value2 = 0.007593534374137497
$ws.Range("H8").Value2 = -0.00598811995247981
$ws.Range("I8").Value2 = 0.009724550898203592
$ws.Range("J8").Value2 = -0.02135822448314467

$ws.Range("C9").Value2 = 0.007447709790839163
$ws.Range("E9").Value2 = 0.02162802251209005
$ws.Range("F9").Value2 = 0.01877556310225241
$ws.Range("H9").Value2 = 0.005200436801747207
$ws.Range("I9").Value2 = 0.02101045204180817
$ws.Range("J9").Value2 = -0.05024317339153084

$ws.Range("C10").Value2 = -0.07582916731666926
$ws.Range("E10").Value2 = -0.03130150120600483
$ws.Range("F10").Value2 = -0.02986485545942184
$ws.Range("H10").Value2 = -0.01797952791811167
$ws.Range("I10").Value2 = -0.02989346357385429
$ws.Range("J10").Value2 = 0.01959309664575326

$ws.Range("C11").Value2 = 0.05106980427921711
$ws.Range("E11").Value2 = -0.06676682706730826
$ws.Range("F11").Value2 = -0.06948910995643984
$ws.Range("H11").Value2 = -0.05326523706094825
$ws.Range("I11").Value2 = -0.06519837279349118
$ws.Range("J11").Value2 = 0.03231384765188664

$ws.Range("C12").Value2 = 0.1043678734714939
$ws.Range("E12").Value2 = 0.5147194028776115
$ws.Range("F12").Value2 = 0.4986034344137377
$ws.Range("H12").Value2 = 0.04393918375673502
$ws.Range("I12").Value2 = 0.5173999255997023
$ws.Range("J12").Value2 = -0.06635186745078865

$ws.Range("C13").Value2 = -0.1727335869343477
$ws.Range("E13").Value2 = -0.02554877419509678
$ws.Range("F13").Value2 = -0.0258624074496298
$ws.Range("H13").Value2 = -0.04002227208908835
$ws.Range("I13").Value2 = -0.02516909267637071
$ws.Range("J13").Value2 = 0.03358693104030525
